# Generate Report for Archive
# - Update the localization status from "Ready for handoff" to
#   "In Translation" everywhere it appears (Overview!E2/F2, and the
#   "Status" column (C2) on both the zh-cn and de-de sheets).
# - Shrink the corresponding status columns that were sized for the
#   old (longer) text so they match the new, narrower content width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- 1. Update the status text -------------------------------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- 2. Resize the columns that were sized for the old text --------------
# ColumnWidth is quantized by the engine to the nearest 1/6 of a character,
# so 12.5 is the value that reliably yields the narrower target width for
# all three affected ranges below.
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
